$d = $word.ActiveDocument

$newParagraphs = @(
    "Let’s check how this works",
    "‘OMG it’s you?!’",
    "“hey hey hey!”"
)

foreach ($text in $newParagraphs) {
    $endRange = $d.Range($d.Content.End, $d.Content.End)
    $endRange.InsertParagraphAfter()
    $insertPos = $d.Content.End - 1
    $paraRange = $d.Range($insertPos, $insertPos)
    $paraRange.InsertAfter($text)
}
